# Auto-generated edit script: refresh market-price derived columns (H-N)
# on the Leve profit sheets, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 12498.131
$ws.Range("I21").Value = 14211.333
$ws.Range("J21").Value = 10629.182
$ws.Range("K21").Value = 14211.333
$ws.Range("L21").Value = 10629.182
$ws.Range("M21").Value = -13743.333
$ws.Range("N21").Value = -11565.182
$ws.Range("H23").Value = 12498.131
$ws.Range("I23").Value = 14211.333
$ws.Range("J23").Value = 10629.182
$ws.Range("K23").Value = 14211.333
$ws.Range("L23").Value = 10629.182
$ws.Range("M23").Value = -13977.333
$ws.Range("N23").Value = -11097.182
$ws.Range("H74").Value = 2800
$ws.Range("I74").Value = 2500
$ws.Range("K74").Value = 2500
$ws.Range("M74").Value = -1564
$ws.Range("H77").Value = 2800
$ws.Range("I77").Value = 2500
$ws.Range("K77").Value = 12500
$ws.Range("M77").Value = -7820
$ws.Range("H100").Value = 1588.1212
$ws.Range("I100").Value = 966.6667
$ws.Range("J100").Value = 2106
$ws.Range("K100").Value = 966.6667
$ws.Range("L100").Value = 2106
$ws.Range("M100").Value = -425.6667
$ws.Range("N100").Value = -3188
$ws.Range("H118").Value = 575.82355
$ws.Range("J118").Value = 850
$ws.Range("L118").Value = 2550
$ws.Range("N118").Value = -5864
$ws.Range("H125").Value = 961.3333
$ws.Range("I125").Value = 800
$ws.Range("K125").Value = 7200
$ws.Range("M125").Value = -4740
$ws.Range("H129").Value = 843.1667
$ws.Range("I129").Value = 226.90909
$ws.Range("J129").Value = 1199.9474
$ws.Range("K129").Value = 680.72727
$ws.Range("L129").Value = 3599.8422
$ws.Range("M129").Value = 4319.27273
$ws.Range("N129").Value = -13599.8422
$ws.Range("H132").Value = 15222313
$ws.Range("I132").Value = 1884016.9
$ws.Range("J132").Value = 71433704
$ws.Range("K132").Value = 5652050.699999999
$ws.Range("L132").Value = 214301112
$ws.Range("M132").Value = -5649520.699999999
$ws.Range("N132").Value = -214306172
$ws.Range("H138").Value = 2534234.8
$ws.Range("I138").Value = 1081.4667
$ws.Range("J138").Value = 4085145
$ws.Range("K138").Value = 3244.4001
$ws.Range("L138").Value = 12255435
$ws.Range("M138").Value = 1895.5999
$ws.Range("N138").Value = -12265715

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 4249.359
$ws.Range("I32").Value = 2875.1016
$ws.Range("J32").Value = 14785.333
$ws.Range("K32").Value = 2875.1016
$ws.Range("L32").Value = 14785.333
$ws.Range("M32").Value = -2588.1016
$ws.Range("N32").Value = -15359.333
$ws.Range("H61").Value = 1073.6522
$ws.Range("J61").Value = 211
$ws.Range("L61").Value = 211
$ws.Range("N61").Value = -635
$ws.Range("H74").Value = 48794.566
$ws.Range("I74").Value = 50913.1
$ws.Range("J74").Value = 34671
$ws.Range("K74").Value = 50913.1
$ws.Range("L74").Value = 34671
$ws.Range("M74").Value = -50039.1
$ws.Range("N74").Value = -36419
$ws.Range("H77").Value = 48794.566
$ws.Range("I77").Value = 50913.1
$ws.Range("J77").Value = 34671
$ws.Range("K77").Value = 254565.5
$ws.Range("L77").Value = 173355
$ws.Range("M77").Value = -250197.5
$ws.Range("N77").Value = -182091
$ws.Range("H122").Value = 1728.2
$ws.Range("I122").Value = 1311.1818
$ws.Range("K122").Value = 3933.5454
$ws.Range("M122").Value = -1483.5454
$ws.Range("H136").Value = 1073.6522
$ws.Range("J136").Value = 211
$ws.Range("L136").Value = 633
$ws.Range("N136").Value = -5733

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1518.238
$ws.Range("I20").Value = 1274.091
$ws.Range("J20").Value = 1786.8
$ws.Range("K20").Value = 1274.091
$ws.Range("L20").Value = 1786.8
$ws.Range("M20").Value = -1027.091
$ws.Range("N20").Value = -2280.8
$ws.Range("H86").Value = 1558.0322
$ws.Range("I86").Value = 1323.2727
$ws.Range("J86").Value = 2131.889
$ws.Range("K86").Value = 1323.2727
$ws.Range("L86").Value = 2131.889
$ws.Range("M86").Value = -200.2727
$ws.Range("N86").Value = -4377.889
$ws.Range("H89").Value = 1558.0322
$ws.Range("I89").Value = 1323.2727
$ws.Range("J89").Value = 2131.889
$ws.Range("K89").Value = 6616.363499999999
$ws.Range("L89").Value = 10659.445
$ws.Range("M89").Value = -1000.363499999999
$ws.Range("N89").Value = -21891.445

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21278146
$ws.Range("I31").Value = 29412694
$ws.Range("J31").Value = 3173.7693
$ws.Range("K31").Value = 29412694
$ws.Range("L31").Value = 3173.7693
$ws.Range("M31").Value = -29412399
$ws.Range("N31").Value = -3763.7693
$ws.Range("H34").Value = 21278146
$ws.Range("I34").Value = 29412694
$ws.Range("J34").Value = 3173.7693
$ws.Range("K34").Value = 29412694
$ws.Range("L34").Value = 3173.7693
$ws.Range("M34").Value = -29412492
$ws.Range("N34").Value = -3577.7693
$ws.Range("H59").Value = 15500
$ws.Range("J59").Value = 15500
$ws.Range("L59").Value = 15500
$ws.Range("N59").Value = -17790
$ws.Range("H94").Value = 978.5
$ws.Range("I94").Value = 1749.5
$ws.Range("J94").Value = 785.75
$ws.Range("K94").Value = 1749.5
$ws.Range("L94").Value = 785.75
$ws.Range("M94").Value = -1298.5
$ws.Range("N94").Value = -1687.75
$ws.Range("H122").Value = 1772
$ws.Range("I122").Value = 1592
$ws.Range("J122").Value = 2132
$ws.Range("K122").Value = 4776
$ws.Range("L122").Value = 6396
$ws.Range("M122").Value = -2326
$ws.Range("N122").Value = -11296
$ws.Range("H132").Value = 2108.718
$ws.Range("I132").Value = 1665.2333
$ws.Range("J132").Value = 3587
$ws.Range("K132").Value = 4995.699900000001
$ws.Range("L132").Value = 10761
$ws.Range("M132").Value = -2465.699900000001
$ws.Range("N132").Value = -15821

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 961.125
$ws.Range("I122").Value = 516.875
$ws.Range("J122").Value = 1849.625
$ws.Range("K122").Value = 4651.875
$ws.Range("L122").Value = 16646.625
$ws.Range("M122").Value = -2201.875
$ws.Range("N122").Value = -21546.625
$ws.Range("H129").Value = 2812.139
$ws.Range("I129").Value = 2690.3076
$ws.Range("J129").Value = 2881
$ws.Range("K129").Value = 8070.9228
$ws.Range("L129").Value = 8643
$ws.Range("M129").Value = -3070.9228
$ws.Range("N129").Value = -18643
$ws.Range("H132").Value = 4680801
$ws.Range("I132").Value = 2085526.4
$ws.Range("J132").Value = 18522266
$ws.Range("K132").Value = 18769737.6
$ws.Range("L132").Value = 166700394
$ws.Range("M132").Value = -18767207.6
$ws.Range("N132").Value = -166705454

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 132.88235
$ws.Range("I2").Value = 54.5
$ws.Range("J2").Value = 202.55556
$ws.Range("K2").Value = 54.5
$ws.Range("L2").Value = 202.55556
$ws.Range("M2").Value = 58.5
$ws.Range("N2").Value = -428.55556
$ws.Range("H113").Value = 1000.5294
$ws.Range("I113").Value = 853
$ws.Range("J113").Value = 1271
$ws.Range("K113").Value = 853
$ws.Range("L113").Value = 1271
$ws.Range("M113").Value = 1317
$ws.Range("N113").Value = -5611

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1018
$ws.Range("I7").Value = 1018
$ws.Range("K7").Value = 1018
$ws.Range("M7").Value = -906
$ws.Range("H46").Value = 947.95
$ws.Range("I46").Value = 1106.125
$ws.Range("J46").Value = 842.5
$ws.Range("K46").Value = 1106.125
$ws.Range("L46").Value = 842.5
$ws.Range("M46").Value = -918.125
$ws.Range("N46").Value = -1218.5
$ws.Range("H55").Value = 491.17392
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 573.5263
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 573.5263
$ws.Range("M55").Value = 73
$ws.Range("N55").Value = -919.5263
$ws.Range("H126").Value = 1018
$ws.Range("I126").Value = 1018
$ws.Range("K126").Value = 3054
$ws.Range("M126").Value = -584

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1012.625
$ws.Range("J81").Value = 650
$ws.Range("L81").Value = 1300
$ws.Range("N81").Value = -3422
$ws.Range("H84").Value = 1012.625
$ws.Range("J84").Value = 650
$ws.Range("L84").Value = 6500
$ws.Range("N84").Value = -17108
